# Item.xlsx edit script
# - Append a new item row (row 9: Equip_Weapon_1 / 开山斧) to Sheet1
# - Grow the XML table (表1) so the new row is included, and refresh the selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Add the new data row under the existing 8 rows of data.
$ws.Range("A9").Value = "Equip_Weapon_1"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "开山斧"
$ws.Range("F9").Value = "开山斧武器"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "50004"
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 10000
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 100

# 2) Grow the bound XML table (表1) to cover the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K9"))

# 3) Move the active selection like the author's last save.
$ws.Range("K13").Select() | Out-Null
